$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '88.015.96'
$ws.Range("E2").Value = '  -3.09%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.107.83'
$ws.Range("E3").Value = '  -2.72%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.28'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '634.85'
$ws.Range("E6").Value = '  +2.66%  '
$ws.Range("E7").Value = '  -1.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.780'
$ws.Range("E8").Value = '  +12.70%  '
$ws.Range("E9").Value = '  +0.07%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.105.07'
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.560'
$ws.Range("E11").Value = '  -3.19%  '
$ws.Range("E12").Value = '  +0.73%  '
$ws.Range("E13").Value = '  -2.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.957.67'
$ws.Range("E15").Value = '  -2.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.681.35'
$ws.Range("E16").Value = '  -2.39%  '
$ws.Range("E17").Value = '  -2.90%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.123.15'
$ws.Range("E18").Value = '  -2.00%  '
$ws.Range("E19").Value = '  +2.89%  '
$ws.Range("E20").Value = '  +15.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.18'
$ws.Range("E21").Value = '  -2.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '420.44'
$ws.Range("E22").Value = '  -5.82%  '
$ws.Range("E23").Value = '  -2.78%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.88'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.42'
$ws.Range("E25").Value = '  +4.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '81.84'
$ws.Range("E26").Value = '  +8.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.41'
$ws.Range("E27").Value = '  -3.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.283.70'
$ws.Range("E28").Value = '  -2.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("E31").Value = '  -8.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.01'
$ws.Range("E32").Value = '  -4.86%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '8.13'
$ws.Range("E33").Value = '  -4.07%  '
$ws.Range("B34").Value = 'Bittensor'
$ws.Range("C34").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '499.43'
$ws.Range("E34").Value = '  -6.63%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.146'
$ws.Range("E35").Value = '  +14.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.86'
$ws.Range("E36").Value = '  -2.22%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").Value = '  +0.17%  '
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.83'
$ws.Range("E38").Value = '  -2.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '22.11'
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("E41").Value = '  +0.32%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.363'
$ws.Range("E43").Value = '  -3.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.84'
$ws.Range("E44").Value = '  -4.47%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.81'
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.133'
$ws.Range("E46").Value = '  +6.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.50'
$ws.Range("E47").Value = '  -1.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0653'
$ws.Range("E48").Value = '  +10.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '161.65'
$ws.Range("E49").Value = '  -6.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.713'
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("E51").Value = '  -5.56%  '
